# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the commit "Created functions to get season record": the sheet
# previously only had team/player statistics through column AC; this adds
# AD:AF with the team's season Wins/Losses/Ties repeated for every player
# row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) should look exactly like the existing headers
# (bold, centered, thin border) -- copy the format from an existing header
# cell, then set the three new labels.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-44: same season record (98 wins, 64 losses, 0 ties) for
# every player on the roster.
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 98
    $ws.Cells.Item($r, 31).Value = 64
    $ws.Cells.Item($r, 32).Value = 0
}
